$d = $word.ActiveDocument

$d.Content.Find.Execute("78+2=80", $true, $false, $false, $false, $false, $true, 1, $false, "27-25=2", 2) | Out-Null
$d.Content.Find.Execute("98-63=35", $true, $false, $false, $false, $false, $true, 1, $false, "31+19=50", 2) | Out-Null
$d.Content.Find.Execute("3+57=60", $true, $false, $false, $false, $false, $true, 1, $false, "51+39=90", 2) | Out-Null
$d.Content.Find.Execute("44+15=59", $true, $false, $false, $false, $false, $true, 1, $false, "3+57=60", 2) | Out-Null
$d.Content.Find.Execute("73-38=35", $true, $false, $false, $false, $false, $true, 1, $false, "62-53=9", 2) | Out-Null
$d.Content.Find.Execute("46+8=54", $true, $false, $false, $false, $false, $true, 1, $false, "67-2=65", 2) | Out-Null
$d.Content.Find.Execute("62+18=80", $true, $false, $false, $false, $false, $true, 1, $false, "33-0=33", 2) | Out-Null
$d.Content.Find.Execute("63-60=3", $true, $false, $false, $false, $false, $true, 1, $false, "85-23=62", 2) | Out-Null
$d.Content.Find.Execute("29+7=36", $true, $false, $false, $false, $false, $true, 1, $false, "3+33=36", 2) | Out-Null
$d.Content.Find.Execute("2+27=29", $true, $false, $false, $false, $false, $true, 1, $false, "58+37=95", 2) | Out-Null
$d.Content.Find.Execute("33+64=97", $true, $false, $false, $false, $false, $true, 1, $false, "5+73=78", 2) | Out-Null
$d.Content.Find.Execute("54-15=39", $true, $false, $false, $false, $false, $true, 1, $false, "18+64=82", 2) | Out-Null
$d.Content.Find.Execute("5+84=89", $true, $false, $false, $false, $false, $true, 1, $false, "4+21=25", 2) | Out-Null
$d.Content.Find.Execute("26+35=61", $true, $false, $false, $false, $false, $true, 1, $false, "59-5=54", 2) | Out-Null
$d.Content.Find.Execute("41+20=61", $true, $false, $false, $false, $false, $true, 1, $false, "58-31=27", 2) | Out-Null
$d.Content.Find.Execute("14+6=20", $true, $false, $false, $false, $false, $true, 1, $false, "55+25=80", 2) | Out-Null
$d.Content.Find.Execute("92+0=92", $true, $false, $false, $false, $false, $true, 1, $false, "89+2=91", 2) | Out-Null
$d.Content.Find.Execute("35+13=48", $true, $false, $false, $false, $false, $true, 1, $false, "97-60=37", 2) | Out-Null
$d.Content.Find.Execute("52-1=51", $true, $false, $false, $false, $false, $true, 1, $false, "54-7=47", 2) | Out-Null
$d.Content.Find.Execute("99-62=37", $true, $false, $false, $false, $false, $true, 1, $false, "17+62=79", 2) | Out-Null
$d.Content.Find.Execute("23+13=36", $true, $false, $false, $false, $false, $true, 1, $false, "42+43=85", 2) | Out-Null
$d.Content.Find.Execute("24-8=16", $true, $false, $false, $false, $false, $true, 1, $false, "91-14=77", 2) | Out-Null
$d.Content.Find.Execute("77-43=34", $true, $false, $false, $false, $false, $true, 1, $false, "58+1=59", 2) | Out-Null
$d.Content.Find.Execute("70+19=89", $true, $false, $false, $false, $false, $true, 1, $false, "46-17=29", 2) | Out-Null
$d.Content.Find.Execute("94-5=89", $true, $false, $false, $false, $false, $true, 1, $false, "15-3=12", 2) | Out-Null
$d.Content.Find.Execute("87-59=28", $true, $false, $false, $false, $false, $true, 1, $false, "86+0=86", 2) | Out-Null
$d.Content.Find.Execute("50+22=72", $true, $false, $false, $false, $false, $true, 1, $false, "15-2=13", 2) | Out-Null
$d.Content.Find.Execute("63+12=75", $true, $false, $false, $false, $false, $true, 1, $false, "77+8=85", 2) | Out-Null
$d.Content.Find.Execute("13+76=89", $true, $false, $false, $false, $false, $true, 1, $false, "86-81=5", 2) | Out-Null
$d.Content.Find.Execute("28-19=9", $true, $false, $false, $false, $false, $true, 1, $false, "59-13=46", 2) | Out-Null
$d.Content.Find.Execute("64-9=55", $true, $false, $false, $false, $false, $true, 1, $false, "84-81=3", 2) | Out-Null
$d.Content.Find.Execute("2+33=35", $true, $false, $false, $false, $false, $true, 1, $false, "48-47=1", 2) | Out-Null
$d.Content.Find.Execute("21-1=20", $true, $false, $false, $false, $false, $true, 1, $false, "55+8=63", 2) | Out-Null
$d.Content.Find.Execute("31+59=90", $true, $false, $false, $false, $false, $true, 1, $false, "83-52=31", 2) | Out-Null
$d.Content.Find.Execute("69-34=35", $true, $false, $false, $false, $false, $true, 1, $false, "16+64=80", 2) | Out-Null
$d.Content.Find.Execute("33-31=2", $true, $false, $false, $false, $false, $true, 1, $false, "44-8=36", 2) | Out-Null
$d.Content.Find.Execute("22+39=61", $true, $false, $false, $false, $false, $true, 1, $false, "96-13=83", 2) | Out-Null
$d.Content.Find.Execute("19+10=29", $true, $false, $false, $false, $false, $true, 1, $false, "1+54=55", 2) | Out-Null
$d.Content.Find.Execute("71-64=7", $true, $false, $false, $false, $false, $true, 1, $false, "95-19=76", 2) | Out-Null
$d.Content.Find.Execute("76-23=53", $true, $false, $false, $false, $false, $true, 1, $false, "5+44=49", 2) | Out-Null
$d.Content.Find.Execute("90-2=88", $true, $false, $false, $false, $false, $true, 1, $false, "3+75=78", 2) | Out-Null
$d.Content.Find.Execute("14+11=25", $true, $false, $false, $false, $false, $true, 1, $false, "40-33=7", 2) | Out-Null
$d.Content.Find.Execute("94-18=76", $true, $false, $false, $false, $false, $true, 1, $false, "74-64=10", 2) | Out-Null
$d.Content.Find.Execute("49+39=88", $true, $false, $false, $false, $false, $true, 1, $false, "15-5=10", 2) | Out-Null
$d.Content.Find.Execute("61+37=98", $true, $false, $false, $false, $false, $true, 1, $false, "70-27=43", 2) | Out-Null
$d.Content.Find.Execute("10+29=39", $true, $false, $false, $false, $false, $true, 1, $false, "68-40=28", 2) | Out-Null
$d.Content.Find.Execute("71+6=77", $true, $false, $false, $false, $false, $true, 1, $false, "81-75=6", 2) | Out-Null
$d.Content.Find.Execute("96-62=34", $true, $false, $false, $false, $false, $true, 1, $false, "10+39=49", 2) | Out-Null
$d.Content.Find.Execute("71-66=5", $true, $false, $false, $false, $false, $true, 1, $false, "60+21=81", 2) | Out-Null
$d.Content.Find.Execute("89-67=22", $true, $false, $false, $false, $false, $true, 1, $false, "14+42=56", 2) | Out-Null
$d.Content.Find.Execute("72-16=56", $true, $false, $false, $false, $false, $true, 1, $false, "27+55=82", 2) | Out-Null
$d.Content.Find.Execute("40-4=36", $true, $false, $false, $false, $false, $true, 1, $false, "84-50=34", 2) | Out-Null
$d.Content.Find.Execute("8+61=69", $true, $false, $false, $false, $false, $true, 1, $false, "75+4=79", 2) | Out-Null
$d.Content.Find.Execute("20+71=91", $true, $false, $false, $false, $false, $true, 1, $false, "58-6=52", 2) | Out-Null
$d.Content.Find.Execute("45-16=29", $true, $false, $false, $false, $false, $true, 1, $false, "25+60=85", 2) | Out-Null
$d.Content.Find.Execute("30+58=88", $true, $false, $false, $false, $false, $true, 1, $false, "90-40=50", 2) | Out-Null
$d.Content.Find.Execute("60-54=6", $true, $false, $false, $false, $false, $true, 1, $false, "53-49=4", 2) | Out-Null
$d.Content.Find.Execute("11+25=36", $true, $false, $false, $false, $false, $true, 1, $false, "12+55=67", 2) | Out-Null
$d.Content.Find.Execute("47+13=60", $true, $false, $false, $false, $false, $true, 1, $false, "11+26=37", 2) | Out-Null
$d.Content.Find.Execute("33+1=34", $true, $false, $false, $false, $false, $true, 1, $false, "16+31=47", 2) | Out-Null
$d.Content.Find.Execute("12+83=95", $true, $false, $false, $false, $false, $true, 1, $false, "83-57=26", 2) | Out-Null
$d.Content.Find.Execute("14+19=33", $true, $false, $false, $false, $false, $true, 1, $false, "76-19=57", 2) | Out-Null
$d.Content.Find.Execute("95-0=95", $true, $false, $false, $false, $false, $true, 1, $false, "82+17=99", 2) | Out-Null
$d.Content.Find.Execute("75-40=35", $true, $false, $false, $false, $false, $true, 1, $false, "50+0=50", 2) | Out-Null
$d.Content.Find.Execute("57-28=29", $true, $false, $false, $false, $false, $true, 1, $false, "16+51=67", 2) | Out-Null
$d.Content.Find.Execute("77-23=54", $true, $false, $false, $false, $false, $true, 1, $false, "8+34=42", 2) | Out-Null
$d.Content.Find.Execute("50-12=38", $true, $false, $false, $false, $false, $true, 1, $false, "39-2=37", 2) | Out-Null
$d.Content.Find.Execute("8+55=63", $true, $false, $false, $false, $false, $true, 1, $false, "66-52=14", 2) | Out-Null
$d.Content.Find.Execute("98-14=84", $true, $false, $false, $false, $false, $true, 1, $false, "90-11=79", 2) | Out-Null
$d.Content.Find.Execute("43+8=51", $true, $false, $false, $false, $false, $true, 1, $false, "18+3=21", 2) | Out-Null
$d.Content.Find.Execute("1+1=2", $true, $false, $false, $false, $false, $true, 1, $false, "66+3=69", 2) | Out-Null
$d.Content.Find.Execute("21+78=99", $true, $false, $false, $false, $false, $true, 1, $false, "41-9=32", 2) | Out-Null
$d.Content.Find.Execute("46-3=43", $true, $false, $false, $false, $false, $true, 1, $false, "49-19=30", 2) | Out-Null
$d.Content.Find.Execute("0+34=34", $true, $false, $false, $false, $false, $true, 1, $false, "80-44=36", 2) | Out-Null
$d.Content.Find.Execute("59-23=36", $true, $false, $false, $false, $false, $true, 1, $false, "90-86=4", 2) | Out-Null
$d.Content.Find.Execute("26+6=32", $true, $false, $false, $false, $false, $true, 1, $false, "63+8=71", 2) | Out-Null
$d.Content.Find.Execute("91-77=14", $true, $false, $false, $false, $false, $true, 1, $false, "73-71=2", 2) | Out-Null
$d.Content.Find.Execute("17+59=76", $true, $false, $false, $false, $false, $true, 1, $false, "65+30=95", 2) | Out-Null
$d.Content.Find.Execute("50-15=35", $true, $false, $false, $false, $false, $true, 1, $false, "95-60=35", 2) | Out-Null
$d.Content.Find.Execute("8+10=18", $true, $false, $false, $false, $false, $true, 1, $false, "18+71=89", 2) | Out-Null
$d.Content.Find.Execute("39-16=23", $true, $false, $false, $false, $false, $true, 1, $false, "77-54=23", 2) | Out-Null
$d.Content.Find.Execute("64-13=51", $true, $false, $false, $false, $false, $true, 1, $false, "95-43=52", 2) | Out-Null
$d.Content.Find.Execute("36+4=40", $true, $false, $false, $false, $false, $true, 1, $false, "6+7=13", 2) | Out-Null
$d.Content.Find.Execute("3+65=68", $true, $false, $false, $false, $false, $true, 1, $false, "79-74=5", 2) | Out-Null
$d.Content.Find.Execute("91-72=19", $true, $false, $false, $false, $false, $true, 1, $false, "86-69=17", 2) | Out-Null
$d.Content.Find.Execute("56-37=19", $true, $false, $false, $false, $false, $true, 1, $false, "69-5=64", 2) | Out-Null
$d.Content.Find.Execute("86-75=11", $true, $false, $false, $false, $false, $true, 1, $false, "88+7=95", 2) | Out-Null
$d.Content.Find.Execute("96-12=84", $true, $false, $false, $false, $false, $true, 1, $false, "96-89=7", 2) | Out-Null
$d.Content.Find.Execute("33-28=5", $true, $false, $false, $false, $false, $true, 1, $false, "89-7=82", 2) | Out-Null
$d.Content.Find.Execute("31-24=7", $true, $false, $false, $false, $false, $true, 1, $false, "32+65=97", 2) | Out-Null
$d.Content.Find.Execute("76-72=4", $true, $false, $false, $false, $false, $true, 1, $false, "60+12=72", 2) | Out-Null
$d.Content.Find.Execute("89-2=87", $true, $false, $false, $false, $false, $true, 1, $false, "52-19=33", 2) | Out-Null
$d.Content.Find.Execute("26-3=23", $true, $false, $false, $false, $false, $true, 1, $false, "40+58=98", 2) | Out-Null
$d.Content.Find.Execute("28+7=35", $true, $false, $false, $false, $false, $true, 1, $false, "93-81=12", 2) | Out-Null
$d.Content.Find.Execute("0+52=52", $true, $false, $false, $false, $false, $true, 1, $false, "79-55=24", 2) | Out-Null
$d.Content.Find.Execute("30-19=11", $true, $false, $false, $false, $false, $true, 1, $false, "66+13=79", 2) | Out-Null
$d.Content.Find.Execute("79-22=57", $true, $false, $false, $false, $false, $true, 1, $false, "88-54=34", 2) | Out-Null
$d.Content.Find.Execute("87-65=22", $true, $false, $false, $false, $false, $true, 1, $false, "19+19=38", 2) | Out-Null
$d.Content.Find.Execute("41+53=94", $true, $false, $false, $false, $false, $true, 1, $false, "29+9=38", 2) | Out-Null
$d.Content.Find.Execute("10+70=80", $true, $false, $false, $false, $false, $true, 1, $false, "62-50=12", 2) | Out-Null
